$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep their original (General) numeric auto-detection
# from converting numeric-looking text (e.g. "1.00", "0.998") into real numbers:
# format the target range as Text before writing, then restore the default
# "Normal" style afterwards so the saved style index matches the original file.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range('D2').Value = '61.153.09'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').Value = '2.643.08'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('D5').Value = '533.01'
$ws.Range('E5').Value = '  +4.49%  '
$ws.Range('D6').Value = '156.57'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = '6.70'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = '0.110'
$ws.Range('E10').Value = '  +6.03%  '
$ws.Range('D11').Value = '0.350'
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '3.102.23'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').Value = '61.147.82'
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').Value = '22.06'
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('E16').Value = '  +4.09%  '
$ws.Range('D17').Value = '2.639.19'
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').Value = '4.80'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '357.29'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').Value = '6.27'
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '61.95'
$ws.Range('E23').Value = '  +2.81%  '
$ws.Range('D24').Value = '0.433'
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  +2.13%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0₃0876'
$ws.Range('E27').Value = '  +4.00%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '7.47'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '6.19'
$ws.Range('E30').Value = '  +7.73%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '19.61'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.63'
$ws.Range('E32').Value = '  +4.28%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '151.31'
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '4.20'
$ws.Range('E34').Value = '  +4.38%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '1.21'
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '0.929'
$ws.Range('E36').Value = '  +10.05%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').Value = '0.888'
$ws.Range('E37').Value = '  +2.88%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.51'
$ws.Range('E38').Value = '  +2.01%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Value = '3.83'
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '297.81'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.646'
$ws.Range('E41').Value = '  +4.42%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.103'
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '5.16'
$ws.Range('E43').Value = '  +5.83%  '
$ws.Range('D44').Value = '0.0567'
$ws.Range('E44').Value = '  +2.29%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '0.997'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '19.95'
$ws.Range('E46').Value = '  +1.15%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0240'
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '19.22'
$ws.Range('E48').Value = '  +7.16%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').Value = '10.35'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '1.87'
$ws.Range('E50').Value = '  +4.88%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.986.66'
$ws.Range('E51').Value = '  -0.48%  '

$editRange.Style = "Normal"
